$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich text cells), preserving other runs ---
# A8 shared string: "Volume 29   Number  44" -> "...45" (chars 21-22 = "44")
$ws.Range("A8").Characters(21,2).Text = "45"

# C9 shared string: "Report Covering the Week  10/31/2022  Through  11/6/2022"
# Replace right-to-left so earlier offsets stay valid:
# "11/6/2022" starts at char 48 (len 9) -> "11/13/2022"
$ws.Range("C9").Characters(48,9).Text = "11/13/2022"
# "10/31/2022" starts at char 27 (len 10) -> "11/7/2022"
$ws.Range("C9").Characters(27,10).Text = "11/7/2022"

# --- Column E width (narrower after new data recalculated best-fit) ---
$ws.Columns("E").ColumnWidth = 7.433768

# --- Cells whose data TYPE flips between number and "no data" text placeholder. ---
# Copy formatting (and, where the destination becomes text, the value too) from a
# stable donor cell elsewhere in the table that already carries the desired style,
# so we reuse the workbooks existing style entries instead of minting new ones.

# C23: 2 -> "0" (no data)
$ws.Range("C14").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C23").PasteSpecial(-4163)

# D23: "0" -> 1
$ws.Range("G14").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("D23").Value = 1

# E23: "***.*" -> -100
$ws.Range("H14").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("E23").Value = -100

# C27: "0" -> 1
$ws.Range("G14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 1

# C28: 1 -> "0" (no data)
$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial(-4163)

# C29: 1 -> "0" (no data)
$ws.Range("C14").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C29").PasteSpecial(-4163)

# F30: "0" -> 1
$ws.Range("G14").Copy()
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("F30").Value = 1

# --- Remaining data table updates (rows 14-30), newly collected crime data ---
# Row 14
$ws.Range("M14").Value = -76.923076923076
$ws.Range("N14").Value = -91.891891891891

# Row 15
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 25
$ws.Range("K15").Value = 44
$ws.Range("M15").Value = 16.129032258064
$ws.Range("N15").Value = -33.333333333333

# Row 16
$ws.Range("C16").Value = 10
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = 42.857142857142
$ws.Range("F16").Value = 24
$ws.Range("G16").Value = 25
$ws.Range("H16").Value = -4
$ws.Range("I16").Value = 283
$ws.Range("J16").Value = 213
$ws.Range("K16").Value = 32.863849765258
$ws.Range("L16").Value = 28.636363636363
$ws.Range("M16").Value = -26.302083333333
$ws.Range("N16").Value = -78.511769172361

# Row 17
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 166.666666666667
$ws.Range("F17").Value = 43
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = 95.454545454545
$ws.Range("I17").Value = 328
$ws.Range("J17").Value = 281
$ws.Range("K17").Value = 16.725978647686
$ws.Range("L17").Value = 18.411552346570
$ws.Range("M17").Value = -2.670623145400
$ws.Range("N17").Value = -55.735492577597

# Row 18
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 26
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = 36.842105263157
$ws.Range("I18").Value = 264
$ws.Range("J18").Value = 232
$ws.Range("K18").Value = 13.793103448275
$ws.Range("L18").Value = -9.897610921501
$ws.Range("M18").Value = -33.501259445843
$ws.Range("N18").Value = -77.280550774526

# Row 19
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = -6.666666666666
$ws.Range("F19").Value = 78
$ws.Range("G19").Value = 53
$ws.Range("H19").Value = 47.169811320754
$ws.Range("I19").Value = 622
$ws.Range("J19").Value = 486
$ws.Range("K19").Value = 27.983539094650
$ws.Range("L19").Value = 32.340425531914
$ws.Range("M19").Value = 133.834586466165
$ws.Range("N19").Value = 27.983539094650

# Row 20
$ws.Range("C20").Value = 5
$ws.Range("E20").Value = 66.666666666666
$ws.Range("F20").Value = 24
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = 71.428571428571
$ws.Range("I20").Value = 218
$ws.Range("J20").Value = 150
$ws.Range("K20").Value = 45.333333333333
$ws.Range("L20").Value = 48.299319727891
$ws.Range("M20").Value = 49.315068493150
$ws.Range("N20").Value = -73.860911270983

# Row 21
$ws.Range("C21").Value = 39
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = 21.875
$ws.Range("F21").Value = 198
$ws.Range("G21").Value = 137
$ws.Range("H21").Value = 44.525547445255
$ws.Range("I21").Value = 1754
$ws.Range("J21").Value = 1393
$ws.Range("K21").Value = 25.915290739411
$ws.Range("L21").Value = 22.657342657342
$ws.Range("M21").Value = 11.43583227446
$ws.Range("N21").Value = -62.124811055927

# Row 22
$ws.Range("C22").Value = 2
$ws.Range("E22").Value = 100
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 33.333333333333
$ws.Range("I22").Value = 13
$ws.Range("J22").Value = 19
$ws.Range("K22").Value = -31.578947368421
$ws.Range("L22").Value = 8.333333333333
$ws.Range("M22").Value = -40.909090909090

# Row 23
$ws.Range("J23").Value = 30
$ws.Range("K23").Value = -16.666666666666

# Row 24
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = -35.714285714285
$ws.Range("F24").Value = 91
$ws.Range("G24").Value = 98
$ws.Range("H24").Value = -7.142857142857
$ws.Range("I24").Value = 928
$ws.Range("J24").Value = 824
$ws.Range("K24").Value = 12.621359223301
$ws.Range("L24").Value = -7.385229540918
$ws.Range("M24").Value = 39.339339339339

# Row 25
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = 20
$ws.Range("F25").Value = 42
$ws.Range("G25").Value = 46
$ws.Range("H25").Value = -8.695652173913
$ws.Range("I25").Value = 519
$ws.Range("J25").Value = 455
$ws.Range("K25").Value = 14.065934065934
$ws.Range("L25").Value = 20.417633410672
$ws.Range("M25").Value = -26.068376068376

# Row 26
$ws.Range("F26").Value = 5
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = 66.666666666666
$ws.Range("I26").Value = 50
$ws.Range("J26").Value = 38
$ws.Range("K26").Value = 31.578947368421
$ws.Range("L26").Value = 72.413793103448

# Row 27
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -42.857142857142
$ws.Range("I27").Value = 59
$ws.Range("J27").Value = 44
$ws.Range("K27").Value = 34.090909090909
$ws.Range("L27").Value = 5.357142857142

# Row 28
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = -33.333333333333
$ws.Range("J28").Value = 37
$ws.Range("K28").Value = -62.162162162162
$ws.Range("M28").Value = -67.441860465116
$ws.Range("N28").Value = -93.269230769230

# Row 29
$ws.Range("E29").Value = -100
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = -33.333333333333
$ws.Range("J29").Value = 27
$ws.Range("K29").Value = -48.148148148148
$ws.Range("M29").Value = -57.575757575757
$ws.Range("N29").Value = -92.746113989637

# Row 30
$ws.Range("I30").Value = 5
$ws.Range("K30").Value = 66.666666666666
$ws.Range("L30").Value = 400

